# ABC Bank Churn Analysis deck - restructure titles/content per new outline.
#
# Helper: replace a single-paragraph shape's text (e.g. a slide title).
# Clearing the paragraph's text before writing the final string keeps the
# interop layer's text-diffing from doing a character-level diff against the
# old text (which can otherwise split the new text into several runs when it
# shares a common substring with what used to be there).
function Set-Title {
    param($shape, $text)
    $tr = $shape.TextFrame.TextRange
    $tr.Paragraphs(1, 1).Text = ""
    $tr.Paragraphs(1, 1).Text = $text
}

# Helper: rewrite a shape's text as N paragraphs, the first of which is left
# blank (an empty <a:p/>), with the remaining paragraphs indented one level
# (IndentLevel = 2 -> a:pPr lvl="1") and holding the given bullet text.
#
# We first fill every paragraph with a short placeholder and only then set
# the final text per-paragraph (instead of joining everything with "`r" up
# front) because the interop layer's text-diffing can otherwise split a
# paragraph's final word into its own run when the new text happens to share
# a trailing substring with pre-existing content.
function Set-ContentBullets {
    param($shape, $bullets)
    $n = $bullets.Count
    $placeholders = @()
    for ($i = 1; $i -le $n; $i++) { $placeholders += "x" }
    $tr = $shape.TextFrame.TextRange
    $tr.Text = [string]::Join("`r", $placeholders)
    for ($i = 2; $i -le $n; $i++) {
        $para = $tr.Paragraphs($i, 1)
        $para.Text = $bullets[$i - 1]
        $para.IndentLevel = 2
    }
    $tr.Paragraphs(1, 1).Text = ""
}

$p = $ppt.ActivePresentation

# --- Slide 1: Title slide ---
$s1 = $p.Slides.Item(1)
Set-Title $s1.Shapes.Item(1) "ABC Bank Customer Churn Analysis"
$sub1 = $s1.Shapes.Item(2)
$subTr = $sub1.TextFrame.TextRange
$subTr.Text = "x`rx`rx"
$subTr.Paragraphs(1, 1).Text = "Executive Summary | Data-Driven Retention Strategy"
$subTr.Paragraphs(2, 1).Text = ""
$subTr.Paragraphs(3, 1).Text = "Charles Walton – Data Analyst Consultant"

# --- Slide 2: Agenda -> Business Objective ---
$s2 = $p.Slides.Item(2)
Set-Title $s2.Shapes.Item(1) "Business Objective"
Set-ContentBullets $s2.Shapes.Item(2) @(
    "",
    "Identify drivers of customer churn across demographics, products, and regions",
    "Quantify churn risk and prioritize high-impact retention opportunities",
    "Provide executive-ready recommendations to improve customer lifetime value"
)

# --- Slide 3: Key Insights Summary -> Key Insight 1 - Demographics & Geography ---
$s3 = $p.Slides.Item(3)
Set-Title $s3.Shapes.Item(1) "Key Insight 1 – Demographics & Geography"
Set-ContentBullets $s3.Shapes.Item(2) @(
    "",
    "Customers aged 50+ show significantly higher churn propensity",
    "Germany exhibits the highest churn rate compared to other regions",
    "Female customers churn at a higher rate than male customers"
)

# --- Slide 4: Recommendation 1 -> Key Insight 2 - Products & Engagement ---
$s4 = $p.Slides.Item(4)
Set-Title $s4.Shapes.Item(1) "Key Insight 2 – Products & Engagement"
Set-ContentBullets $s4.Shapes.Item(2) @(
    "",
    "Customers with 3–4 products churn more than those with 1–2 products",
    "Inactive members are significantly more likely to exit",
    "Diamond cardholders show the highest churn rate among card types"
)

# --- Slide 5: Recommendation 2 -> Key Insight 3 - Financial Indicators ---
$s5 = $p.Slides.Item(5)
Set-Title $s5.Shapes.Item(1) "Key Insight 3 – Financial Indicators"
Set-ContentBullets $s5.Shapes.Item(2) @(
    "",
    "High-balance customers are more likely to churn than low-balance customers",
    "Lower credit scores correlate with higher churn risk",
    "Tenure reduces churn slightly, but risk persists across all years"
)

# --- Slide 6: Recommendation 3 -> Strategic Recommendations ---
$s6 = $p.Slides.Item(6)
Set-Title $s6.Shapes.Item(1) "Strategic Recommendations"
Set-ContentBullets $s6.Shapes.Item(2) @(
    "",
    "Deploy targeted retention programs for high-balance and senior customers",
    "Redesign Diamond card benefits to improve perceived value",
    "Introduce engagement incentives for inactive members",
    "Launch region-specific interventions for Germany"
)

# --- Slide 7: Recommendation 4 -> Expected Business Impact ---
$s7 = $p.Slides.Item(7)
Set-Title $s7.Shapes.Item(1) "Expected Business Impact"
Set-ContentBullets $s7.Shapes.Item(2) @(
    "",
    "3–5% reduction in annual churn rate",
    "Increased customer lifetime value through improved retention",
    "Improved marketing ROI via targeted interventions",
    "Stronger customer satisfaction and loyalty metrics"
)

# --- Slide 8: Recommendation 5 -> KPI - Overall Churn Rate ---
$s8 = $p.Slides.Item(8)
Set-Title $s8.Shapes.Item(1) "KPI – Overall Churn Rate"
Set-ContentBullets $s8.Shapes.Item(2) @(
    "",
    "Tracks percentage of customers exiting the bank",
    "Primary indicator of retention performance",
    "Used to measure success of churn reduction initiatives"
)

# --- Slide 9: Expected Business Impact -> KPI - Churn by Geography ---
$s9 = $p.Slides.Item(9)
Set-Title $s9.Shapes.Item(1) "KPI – Churn by Geography"
Set-ContentBullets $s9.Shapes.Item(2) @(
    "",
    "Highlights regional risk concentrations",
    "Supports location-specific retention strategies",
    "Germany identified as highest-risk region"
)

# --- Slide 10: Conclusion -> KPI - Product Holding vs Churn ---
$s10 = $p.Slides.Item(10)
Set-Title $s10.Shapes.Item(1) "KPI – Product Holding vs Churn"
Set-ContentBullets $s10.Shapes.Item(2) @(
    "",
    "Measures churn across product ownership levels",
    "Identifies complexity and service friction risks",
    "Supports product simplification and bundling strategies"
)
